$wb = $excel.ActiveWorkbook

# ===== Sheet 1: Overview =====
$ws1 = $wb.Worksheets.Item(1)
$ws1.Hyperlinks.Delete()
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"
$ws1.Range("D2").Value = "2016-03-23 09:42:30"
$ws1.Range("B3").Value = "Handed back: in sync with en-US"
$ws1.Range("C3").Value = "Handed back: in sync with en-US"
$ws1.Range("D3").Value = "2016-03-23 09:40:54"
$ws1.Range("B4").Value = "Handed back: in sync with en-US"
$ws1.Range("C4").Value = "Handed back: in sync with en-US"
$ws1.Range("D4").Value = "2016-03-23 09:42:30"
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/39ae309b41b3e3746f5a5f4a0f42154000ea52b2/e2e/d06b4c91-c8cc-40fd-a5db-02702911b3f4.md", "", "", "d06b4c91-c8cc-40fd-a5db-02702911b3f4.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2e994d3b14fe2273d58c279b4be6f09c3c669583/e2e/d2f80547-3b07-445f-ae15-9c500b9db91d.md", "", "", "d2f80547-3b07-445f-ae15-9c500b9db91d.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/39ae309b41b3e3746f5a5f4a0f42154000ea52b2/e2e/d8ecfe42-2a9b-40e1-8dda-6d704190ec29.md", "", "", "d8ecfe42-2a9b-40e1-8dda-6d704190ec29.md") | Out-Null

# ===== Sheet 2: zh-cn =====
$ws2 = $wb.Worksheets.Item(2)
$ws2.Hyperlinks.Delete()
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("E2").Value = "2016-03-23 09:42:21"
$ws2.Range("H2").Value = "2016-03-23 09:43:03"
$ws2.Range("J2").Value = "Include"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("E3").Value = "2016-03-23 09:40:46"
$ws2.Range("H3").Value = "2016-03-23 09:41:27"
$ws2.Range("J3").Value = "Include"
$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = "Handed back: in sync with en-US"
$ws2.Range("E4").Value = "2016-03-23 09:42:21"
$ws2.Range("H4").Value = "2016-03-23 09:43:03"
$ws2.Range("J4").Value = "Include"
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/39ae309b41b3e3746f5a5f4a0f42154000ea52b2/e2e/d06b4c91-c8cc-40fd-a5db-02702911b3f4.md", "", "", "d06b4c91-c8cc-40fd-a5db-02702911b3f4.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aa409212898c31b45baa723d3982866b5ffd0f8d/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/d06b4c91-c8cc-40fd-a5db-02702911b3f4.0e19f34147afd4266c721ec546cd80a3c6b34dcc.zh-cn.xlf", "", "", "d06b4c91-c8cc-40fd-a5db-02702911b3f4.0e19f34147afd4266c721ec546cd80a3c6b34dcc.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/42491e7fce833b769784bf8832fcdbdc82c0bed1/e2e/d06b4c91-c8cc-40fd-a5db-02702911b3f4.md", "", "", "d06b4c91-c8cc-40fd-a5db-02702911b3f4.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/eae2d5a7020669b29069f5c55db1db9b7d3796d1/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/d06b4c91-c8cc-40fd-a5db-02702911b3f4.0e19f34147afd4266c721ec546cd80a3c6b34dcc.zh-cn.xlf", "", "", "d06b4c91-c8cc-40fd-a5db-02702911b3f4.0e19f34147afd4266c721ec546cd80a3c6b34dcc.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2e994d3b14fe2273d58c279b4be6f09c3c669583/e2e/d2f80547-3b07-445f-ae15-9c500b9db91d.md", "", "", "d2f80547-3b07-445f-ae15-9c500b9db91d.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/835d65a4f09081a87ee5b419ab3a958990ac36af/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.zh-cn.xlf", "", "", "d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/42491e7fce833b769784bf8832fcdbdc82c0bed1/e2e/d2f80547-3b07-445f-ae15-9c500b9db91d.md", "", "", "d2f80547-3b07-445f-ae15-9c500b9db91d.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/eae2d5a7020669b29069f5c55db1db9b7d3796d1/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.zh-cn.xlf", "", "", "d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/39ae309b41b3e3746f5a5f4a0f42154000ea52b2/e2e/d8ecfe42-2a9b-40e1-8dda-6d704190ec29.md", "", "", "d8ecfe42-2a9b-40e1-8dda-6d704190ec29.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aa409212898c31b45baa723d3982866b5ffd0f8d/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/d06b4c91-c8cc-40fd-a5db-02702911b3f4.0e19f34147afd4266c721ec546cd80a3c6b34dcc.zh-cn.xlf", "", "", "d06b4c91-c8cc-40fd-a5db-02702911b3f4.0e19f34147afd4266c721ec546cd80a3c6b34dcc.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/42491e7fce833b769784bf8832fcdbdc82c0bed1/e2e/d06b4c91-c8cc-40fd-a5db-02702911b3f4.md", "", "", "d06b4c91-c8cc-40fd-a5db-02702911b3f4.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/eae2d5a7020669b29069f5c55db1db9b7d3796d1/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/d06b4c91-c8cc-40fd-a5db-02702911b3f4.0e19f34147afd4266c721ec546cd80a3c6b34dcc.zh-cn.xlf", "", "", "d06b4c91-c8cc-40fd-a5db-02702911b3f4.0e19f34147afd4266c721ec546cd80a3c6b34dcc.zh-cn.xlf") | Out-Null

# ===== Sheet 3: de-de =====
$ws3 = $wb.Worksheets.Item(3)
$ws3.Hyperlinks.Delete()
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("E2").Value = "2016-03-23 09:42:30"
$ws3.Range("H2").Value = "2016-03-23 09:43:20"
$ws3.Range("J2").Value = "Include"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("E3").Value = "2016-03-23 09:40:54"
$ws3.Range("H3").Value = "2016-03-23 09:41:42"
$ws3.Range("J3").Value = "Include"
$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = "Handed back: in sync with en-US"
$ws3.Range("E4").Value = "2016-03-23 09:42:30"
$ws3.Range("H4").Value = "2016-03-23 09:43:20"
$ws3.Range("J4").Value = "Include"
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/39ae309b41b3e3746f5a5f4a0f42154000ea52b2/e2e/d06b4c91-c8cc-40fd-a5db-02702911b3f4.md", "", "", "d06b4c91-c8cc-40fd-a5db-02702911b3f4.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a540ccfcc56464f10c0e93ae52721118e781513c/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/d06b4c91-c8cc-40fd-a5db-02702911b3f4.0e19f34147afd4266c721ec546cd80a3c6b34dcc.de-de.xlf", "", "", "d06b4c91-c8cc-40fd-a5db-02702911b3f4.0e19f34147afd4266c721ec546cd80a3c6b34dcc.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/4f60ca16efcb9e19e50b2d110ed0860df54a1b6b/e2e/d06b4c91-c8cc-40fd-a5db-02702911b3f4.md", "", "", "d06b4c91-c8cc-40fd-a5db-02702911b3f4.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/29bef86d852763c7d4224420ef1223a2b3983acc/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/d06b4c91-c8cc-40fd-a5db-02702911b3f4.0e19f34147afd4266c721ec546cd80a3c6b34dcc.de-de.xlf", "", "", "d06b4c91-c8cc-40fd-a5db-02702911b3f4.0e19f34147afd4266c721ec546cd80a3c6b34dcc.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2e994d3b14fe2273d58c279b4be6f09c3c669583/e2e/d2f80547-3b07-445f-ae15-9c500b9db91d.md", "", "", "d2f80547-3b07-445f-ae15-9c500b9db91d.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7cbf60a888a796051527726f6649d50562368f2a/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.de-de.xlf", "", "", "d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/4f60ca16efcb9e19e50b2d110ed0860df54a1b6b/e2e/d2f80547-3b07-445f-ae15-9c500b9db91d.md", "", "", "d2f80547-3b07-445f-ae15-9c500b9db91d.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/29bef86d852763c7d4224420ef1223a2b3983acc/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.de-de.xlf", "", "", "d2f80547-3b07-445f-ae15-9c500b9db91d.7383fb7fb7f6acc8184223d1325040eae2505798.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/39ae309b41b3e3746f5a5f4a0f42154000ea52b2/e2e/d8ecfe42-2a9b-40e1-8dda-6d704190ec29.md", "", "", "d8ecfe42-2a9b-40e1-8dda-6d704190ec29.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a540ccfcc56464f10c0e93ae52721118e781513c/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/d06b4c91-c8cc-40fd-a5db-02702911b3f4.0e19f34147afd4266c721ec546cd80a3c6b34dcc.de-de.xlf", "", "", "d06b4c91-c8cc-40fd-a5db-02702911b3f4.0e19f34147afd4266c721ec546cd80a3c6b34dcc.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F4"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/4f60ca16efcb9e19e50b2d110ed0860df54a1b6b/e2e/d06b4c91-c8cc-40fd-a5db-02702911b3f4.md", "", "", "d06b4c91-c8cc-40fd-a5db-02702911b3f4.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/29bef86d852763c7d4224420ef1223a2b3983acc/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/d06b4c91-c8cc-40fd-a5db-02702911b3f4.0e19f34147afd4266c721ec546cd80a3c6b34dcc.de-de.xlf", "", "", "d06b4c91-c8cc-40fd-a5db-02702911b3f4.0e19f34147afd4266c721ec546cd80a3c6b34dcc.de-de.xlf") | Out-Null
